# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (col I) and DialogAct (col J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;   Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 20;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 24;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 30;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 35;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 39;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 48;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 51;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 53;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 56;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 60;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 64;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 70;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 74;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 76;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 90;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 91;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 92;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 98;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 103; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 105; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 110; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 121; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 125; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 135; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 143; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 144; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 156; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
